$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 3543.5334
$ws.Cells.Item(28, 9).Value = 1001.2222
$ws.Cells.Item(28, 11).Value = 1001.2222
$ws.Cells.Item(28, 13).Value = -516.2222

$ws.Cells.Item(40, 8).Value = 5767.273
$ws.Cells.Item(40, 9).Value = 4206.143
$ws.Cells.Item(40, 11).Value = 4206.143
$ws.Cells.Item(40, 13).Value = -4031.143

$ws.Cells.Item(64, 8).Value = 27783966
$ws.Cells.Item(64, 9).Value = 6353.1724
$ws.Cells.Item(64, 11).Value = 6353.1724
$ws.Cells.Item(64, 13).Value = -6105.1724

$ws.Cells.Item(67, 8).Value = 27783966
$ws.Cells.Item(67, 9).Value = 6353.1724
$ws.Cells.Item(67, 11).Value = 6353.1724
$ws.Cells.Item(67, 13).Value = -5495.1724

$ws.Cells.Item(107, 8).Value = 409.6316
$ws.Cells.Item(107, 9).Value = 273.625
$ws.Cells.Item(107, 10).Value = 1135
$ws.Cells.Item(107, 11).Value = 273.625
$ws.Cells.Item(107, 12).Value = 1135
$ws.Cells.Item(107, 13).Value = 1646.375
$ws.Cells.Item(107, 14).Value = -4975

$ws.Cells.Item(132, 8).Value = 4792.4546
$ws.Cells.Item(132, 9).Value = 4857.4443
$ws.Cells.Item(132, 11).Value = 14572.3329
$ws.Cells.Item(132, 13).Value = -12042.3329

$ws.Cells.Item(133, 8).Value = 105995.5
$ws.Cells.Item(133, 10).Value = 105995.5
$ws.Cells.Item(133, 12).Value = 105995.5
$ws.Cells.Item(133, 14).Value = -116115.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4166.8687
$ws.Cells.Item(32, 9).Value = 3576.3225
$ws.Cells.Item(32, 11).Value = 3576.3225
$ws.Cells.Item(32, 13).Value = -3289.3225

$ws.Cells.Item(45, 8).Value = 35652.77
$ws.Cells.Item(45, 9).Value = 54491.25
$ws.Cells.Item(45, 11).Value = 54491.25
$ws.Cells.Item(45, 13).Value = -54114.25

$ws.Cells.Item(61, 8).Value = 3323.9285
$ws.Cells.Item(61, 9).Value = 3271.923
$ws.Cells.Item(61, 11).Value = 3271.923
$ws.Cells.Item(61, 13).Value = -3059.923

$ws.Cells.Item(74, 8).Value = 224762.2
$ws.Cells.Item(74, 9).Value = 618399.5600000001
$ws.Cells.Item(74, 11).Value = 618399.5600000001
$ws.Cells.Item(74, 13).Value = -617525.5600000001

$ws.Cells.Item(77, 8).Value = 224762.2
$ws.Cells.Item(77, 9).Value = 618399.5600000001
$ws.Cells.Item(77, 11).Value = 3091997.8
$ws.Cells.Item(77, 13).Value = -3087629.8

$ws.Cells.Item(110, 8).Value = 4417.3
$ws.Cells.Item(110, 9).Value = 1732.6666
$ws.Cells.Item(110, 10).Value = 5567.857
$ws.Cells.Item(110, 11).Value = 1732.6666
$ws.Cells.Item(110, 12).Value = 5567.857
$ws.Cells.Item(110, 13).Value = 312.3334
$ws.Cells.Item(110, 14).Value = -9657.857

$ws.Cells.Item(122, 8).Value = 4206.25
$ws.Cells.Item(122, 9).Value = 1595.7646
$ws.Cells.Item(122, 11).Value = 4787.293799999999
$ws.Cells.Item(122, 13).Value = -2337.293799999999

$ws.Cells.Item(132, 8).Value = 2241.7666
$ws.Cells.Item(132, 9).Value = 1337.4348
$ws.Cells.Item(132, 11).Value = 4012.3044
$ws.Cells.Item(132, 13).Value = -1482.3044

$ws.Cells.Item(136, 8).Value = 3323.9285
$ws.Cells.Item(136, 9).Value = 3271.923
$ws.Cells.Item(136, 11).Value = 9815.769
$ws.Cells.Item(136, 13).Value = -7265.769

$ws.Cells.Item(141, 8).Value = 107499.4
$ws.Cells.Item(141, 10).Value = 107499.4
$ws.Cells.Item(141, 12).Value = 107499.4
$ws.Cells.Item(141, 14).Value = -117859.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 940
$ws.Cells.Item(80, 10).Value = 964.6667
$ws.Cells.Item(80, 12).Value = 964.6667
$ws.Cells.Item(80, 14).Value = -2960.6667

$ws.Cells.Item(83, 8).Value = 940
$ws.Cells.Item(83, 10).Value = 964.6667
$ws.Cells.Item(83, 12).Value = 4823.3335
$ws.Cells.Item(83, 14).Value = -14807.3335

$ws.Cells.Item(86, 8).Value = 2160.5
$ws.Cells.Item(86, 9).Value = 1866.8334
$ws.Cells.Item(86, 10).Value = 2747.8333
$ws.Cells.Item(86, 11).Value = 1866.8334
$ws.Cells.Item(86, 12).Value = 2747.8333
$ws.Cells.Item(86, 13).Value = -743.8334
$ws.Cells.Item(86, 14).Value = -4993.8333

$ws.Cells.Item(89, 8).Value = 2160.5
$ws.Cells.Item(89, 9).Value = 1866.8334
$ws.Cells.Item(89, 10).Value = 2747.8333
$ws.Cells.Item(89, 11).Value = 9334.166999999999
$ws.Cells.Item(89, 12).Value = 13739.1665
$ws.Cells.Item(89, 13).Value = -3718.166999999999
$ws.Cells.Item(89, 14).Value = -24971.1665

$ws.Cells.Item(135, 8).Value = 75995.8
$ws.Cells.Item(135, 10).Value = 75995.8
$ws.Cells.Item(135, 12).Value = 75995.8
$ws.Cells.Item(135, 14).Value = -86135.8

$ws.Cells.Item(138, 8).Value = 34237.383
$ws.Cells.Item(138, 10).Value = 34237.383
$ws.Cells.Item(138, 12).Value = 34237.383
$ws.Cells.Item(138, 14).Value = -44517.383

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4277.0605
$ws.Cells.Item(31, 9).Value = 2925.72
$ws.Cells.Item(31, 11).Value = 2925.72
$ws.Cells.Item(31, 13).Value = -2630.72

$ws.Cells.Item(34, 8).Value = 4277.0605
$ws.Cells.Item(34, 9).Value = 2925.72
$ws.Cells.Item(34, 11).Value = 2925.72
$ws.Cells.Item(34, 13).Value = -2723.72

$ws.Cells.Item(58, 8).Value = 2438.3333
$ws.Cells.Item(58, 9).Value = 1372.2858
$ws.Cells.Item(58, 11).Value = 1372.2858
$ws.Cells.Item(58, 13).Value = -1169.2858

$ws.Cells.Item(96, 8).Value = 6968.5713
$ws.Cells.Item(96, 10).Value = 6968.5713
$ws.Cells.Item(96, 12).Value = 6968.5713
$ws.Cells.Item(96, 14).Value = -12460.5713

$ws.Cells.Item(105, 8).Value = 975.1579
$ws.Cells.Item(105, 9).Value = 610.8570999999999
$ws.Cells.Item(105, 11).Value = 610.8570999999999
$ws.Cells.Item(105, 13).Value = 1136.1429

$ws.Cells.Item(134, 8).Value = 3063.2856
$ws.Cells.Item(134, 9).Value = 2901.7058
$ws.Cells.Item(134, 11).Value = 8705.117400000001
$ws.Cells.Item(134, 13).Value = -6170.117400000001

$ws.Cells.Item(136, 8).Value = 2438.3333
$ws.Cells.Item(136, 9).Value = 1372.2858
$ws.Cells.Item(136, 11).Value = 4116.857400000001
$ws.Cells.Item(136, 13).Value = -1566.857400000001

$ws.Cells.Item(141, 8).Value = 655571.3
$ws.Cells.Item(141, 10).Value = 655571.3
$ws.Cells.Item(141, 12).Value = 655571.3
$ws.Cells.Item(141, 14).Value = -665931.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 7634.5
$ws.Cells.Item(87, 9).Value = 269
$ws.Cells.Item(87, 11).Value = 807
$ws.Cells.Item(87, 13).Value = 441

$ws.Cells.Item(90, 8).Value = 7634.5
$ws.Cells.Item(90, 9).Value = 269
$ws.Cells.Item(90, 11).Value = 2421
$ws.Cells.Item(90, 13).Value = 3819

$ws.Cells.Item(116, 8).Value = 131480
$ws.Cells.Item(116, 9).Value = 300398.66
$ws.Cells.Item(116, 10).Value = 4791
$ws.Cells.Item(116, 11).Value = 901195.98
$ws.Cells.Item(116, 12).Value = 14373
$ws.Cells.Item(116, 13).Value = -897753.98
$ws.Cells.Item(116, 14).Value = -21257

$ws.Cells.Item(129, 8).Value = 2998
$ws.Cells.Item(129, 9).Value = 2630
$ws.Cells.Item(129, 10).Value = 3090
$ws.Cells.Item(129, 11).Value = 7890
$ws.Cells.Item(129, 12).Value = 9270
$ws.Cells.Item(129, 13).Value = -2890
$ws.Cells.Item(129, 14).Value = -19270

$ws.Cells.Item(140, 8).Value = 14330.75
$ws.Cells.Item(140, 9).Value = 11208
$ws.Cells.Item(140, 11).Value = 33624
$ws.Cells.Item(140, 13).Value = -28444

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(31, 8).Value = 477
$ws.Cells.Item(31, 9).Value = 477
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 477
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = -185
$ws.Cells.Item(31, 14).ClearContents()

$ws.Cells.Item(37, 8).Value = 477
$ws.Cells.Item(37, 9).Value = 477
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 477
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 13).Value = -200
$ws.Cells.Item(37, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 3954.2917
$ws.Cells.Item(122, 9).Value = 2313.6924
$ws.Cells.Item(122, 11).Value = 6941.0772
$ws.Cells.Item(122, 13).Value = -4491.0772

$ws.Cells.Item(132, 8).Value = 3065
$ws.Cells.Item(132, 9).Value = 2723.8845
$ws.Cells.Item(132, 10).Value = 7499.5
$ws.Cells.Item(132, 11).Value = 8171.6535
$ws.Cells.Item(132, 12).Value = 22498.5
$ws.Cells.Item(132, 13).Value = -5641.6535
$ws.Cells.Item(132, 14).Value = -27558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2194.7058
$ws.Cells.Item(82, 10).Value = 2075
$ws.Cells.Item(82, 12).Value = 2075
$ws.Cells.Item(82, 14).Value = -2797

$ws.Cells.Item(85, 8).Value = 2194.7058
$ws.Cells.Item(85, 10).Value = 2075
$ws.Cells.Item(85, 12).Value = 2075
$ws.Cells.Item(85, 14).Value = -4571

$ws.Cells.Item(122, 8).Value = 5803.28
$ws.Cells.Item(122, 9).Value = 5004.048
$ws.Cells.Item(122, 11).Value = 15012.144
$ws.Cells.Item(122, 13).Value = -12562.144

$ws.Cells.Item(132, 8).Value = 3968.359
$ws.Cells.Item(132, 9).Value = 2737.1482
$ws.Cells.Item(132, 11).Value = 8211.444600000001
$ws.Cells.Item(132, 13).Value = -5681.444600000001

$ws.Cells.Item(136, 8).Value = 4271.9546
$ws.Cells.Item(136, 9).Value = 4544.364
$ws.Cells.Item(136, 11).Value = 13633.092
$ws.Cells.Item(136, 13).Value = -11083.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(104, 8).Value = 25411
$ws.Cells.Item(104, 10).Value = 25411
$ws.Cells.Item(104, 12).Value = 25411
$ws.Cells.Item(104, 14).Value = -32399

$ws.Cells.Item(122, 8).Value = 11906018
$ws.Cells.Item(122, 9).Value = 1252
$ws.Cells.Item(122, 11).Value = 3756
$ws.Cells.Item(122, 13).Value = -1306

$ws.Cells.Item(126, 8).Value = 4999.3335
$ws.Cells.Item(126, 9).Value = 3999
$ws.Cells.Item(126, 11).Value = 11997
$ws.Cells.Item(126, 13).Value = -9527

$ws.Cells.Item(132, 8).Value = 4292.7915
$ws.Cells.Item(132, 9).Value = 3626.75
$ws.Cells.Item(132, 11).Value = 10880.25
$ws.Cells.Item(132, 13).Value = -8350.25

$ws.Cells.Item(136, 8).Value = 1407
$ws.Cells.Item(136, 10).Value = 1649.5
$ws.Cells.Item(136, 12).Value = 4948.5
$ws.Cells.Item(136, 14).Value = -10048.5
